$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the member list table
$ws.Range("A1").Value = "Prénoms"
$ws.Range("B1").Value = "Nom"

$ws.Range("A2").Value = "Khadidiatou"
$ws.Range("B2").Value = "Coulibaly"

$ws.Range("A3").Value = "Tamsir"
$ws.Range("B3").Value = "Ndong"

$ws.Range("A4").Value = "Samba"
$ws.Range("B4").Value = "Dieng"

$ws.Range("A5").Value = "Jeanne De La Flèche"
$ws.Range("B5").Value = "Onanena Amana"

# Column A autofits to its widest content (bestFit), sized to fit "Jeanne De La Flèche"
$ws.Columns.Item(1).ColumnWidth = 16.6

# Selection ends up on M9 after data entry
$ws.Range("M9").Select() | Out-Null
